# IE_settings.xlsx update
# - Fix misaligned ALANG coefficients on the RSE sheet (rolling-process
#   outputs that produce more than one product need a split coefficient
#   smaller than 1 instead of the old placeholder 1/10 values).
# - Bump the R2M "date" entries forward two days (20200426 -> 20200428).
# - Leave the workbook with R2M as the active/selected sheet.

$wb = $excel.ActiveWorkbook

# --- R2M sheet: refresh the date stamps -------------------------------
$r2m = $wb.Worksheets.Item("R2M")
$r2m.Range("B2").Value = 20200428
$r2m.Range("B4").Value = 20200428

# --- RSE sheet: correct the split coefficients -------------------------
$rse = $wb.Worksheets.Item("RSE")

# Supply
$rse.Range("B2").Value = 0.2
$rse.Range("C2").Value = 1

# Use
$rse.Range("B3").Value = 0.5
$rse.Range("C3").Value = 2
$rse.Range("D3").Value = 100

# Extraction
$rse.Range("B5").Value = 0.1
$rse.Range("C5").Value = 0.3

# EolScrap
$rse.Range("B6").Value = 0.2
$rse.Range("C6").Value = 1

# Zero
$rse.Range("D9").Value = 1

# --- Selections / active sheet -----------------------------------------
# RSE keeps a selection at D3 but is no longer the tab in focus.
$rse.Range("D3").Select() | Out-Null

# R2M becomes the selected/active sheet, with the cursor parked at B5.
$r2m.Activate() | Out-Null
$r2m.Range("B5").Select() | Out-Null
